$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($ws, $ref, $val) {
    $cell = $ws.Range($ref)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

$ws.Range("D2").Value = '71.716.36'
$ws.Range("E2").Value = '  +3.37%  '
$ws.Range("D3").Value = '3.689.25'
$ws.Range("E3").Value = '  +8.48%  '
$ws.Range("E4").Value = '  +0.07%  '
Set-TextValue $ws "D5" '589.14'
$ws.Range("E5").Value = '  +1.34%  '
Set-TextValue $ws "D6" '180.34'
$ws.Range("E6").Value = '  +0.61%  '
$ws.Range("D7").Value = '3.679.85'
$ws.Range("E7").Value = '  +8.46%  '
Set-TextValue $ws "D8" '0.622'
$ws.Range("E9").Value = '  +0.00%  '
$ws.Range("E10").Value = '  +1.38%  '
Set-TextValue $ws "D11" '0.614'
$ws.Range("E11").Value = '  +4.62%  '
Set-TextValue $ws "D12" '50.02'
$ws.Range("E12").Value = '  +3.24%  '
Set-TextValue $ws "D13" '0.0000288'
$ws.Range("E13").Value = '  +1.35%  '
$ws.Range("D14").Value = '4.284.08'
$ws.Range("E14").Value = '  +8.53%  '
Set-TextValue $ws "D15" '685.86'
$ws.Range("E15").Value = '  +0.64%  '
Set-TextValue $ws "D16" '9.01'
$ws.Range("E16").Value = '  +4.70%  '
$ws.Range("B17").Value = 'WrappedEther'
$ws.Range("C17").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D17").Value = '3.693.07'
$ws.Range("E17").Value = '  +8.67%  '
$ws.Range("B18").Value = 'WrappedBTC'
$ws.Range("C18").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D18").Value = '71.786.29'
$ws.Range("E18").Value = '  +3.34%  '
$ws.Range("E19").Value = '  +2.18%  '
Set-TextValue $ws "D20" '18.13'
$ws.Range("E20").Value = '  +2.46%  '
Set-TextValue $ws "D21" '11.68'
$ws.Range("E21").Value = '  +3.40%  '
Set-TextValue $ws "D22" '0.941'
$ws.Range("E22").Value = '  +3.40%  '
$ws.Range("E23").Value = '  +17.44%  '
Set-TextValue $ws "D24" '17.84'
$ws.Range("E24").Value = '  +4.39%  '
Set-TextValue $ws "D25" '104.09'
$ws.Range("E25").Value = '  +2.87%  '
$ws.Range("E26").Value = '  +3.66%  '
$ws.Range("E27").Value = '  +5.37%  '
Set-TextValue $ws "D28" '10.22'
$ws.Range("E28").Value = '  +4.89%  '
Set-TextValue $ws "D29" '35.55'
$ws.Range("E29").Value = '  +5.93%  '
Set-TextValue $ws "D30" '9.24'
$ws.Range("E30").Value = '  +5.44%  '
Set-TextValue $ws "D31" '7.37'
$ws.Range("E31").Value = '  +6.60%  '
$ws.Range("E32").Value = '  +12.63%  '
Set-TextValue $ws "D33" '574.11'
$ws.Range("E34").Value = '  +2.54%  '
Set-TextValue $ws "D35" '0.110'
$ws.Range("E35").Value = '  +3.70%  '
Set-TextValue $ws "D36" '59.51'
$ws.Range("E36").Value = '  +2.54%  '
$ws.Range("D37").Value = '3.796.25'
$ws.Range("E37").Value = '  +5.19%  '
$ws.Range("E38").Value = '  -0.07%  '
$ws.Range("E39").Value = '  +4.10%  '
$ws.Range("D40").Value = '0.0₃0780'
$ws.Range("E40").Value = '  +4.04%  '
Set-TextValue $ws "D41" '35.50'
$ws.Range("E41").Value = '  +0.52%  '
$ws.Range("E42").Value = '  +5.84%  '
Set-TextValue $ws "D43" '0.0463'
$ws.Range("E43").Value = '  +8.54%  '
Set-TextValue $ws "D44" '2.80'
$ws.Range("E44").Value = '  +3.50%  '
Set-TextValue $ws "D45" '0.353'
$ws.Range("E45").Value = '  +4.94%  '
Set-TextValue $ws "D46" '2.89'
$ws.Range("E46").Value = '  +7.84%  '
$ws.Range("E47").Value = '  +0.74%  '
$ws.Range("E48").Value = '  +4.13%  '
$ws.Range("E49").Value = '  +2.46%  '
Set-TextValue $ws "D50" '0.998'
$ws.Range("E50").Value = '  -0.18%  '
Set-TextValue $ws "D51" '134.61'
$ws.Range("E51").Value = '  +2.69%  '
